# Daily attendance processing - 2026-01-18 19:31:38
#
# For every row in the "Recorded By" column (G), when the recorded-by
# list is "dnasr281@gmail.com, System", reorder it to "System, dnasr281@gmail.com".
# All other values in that column (e.g. plain "dnasr281@gmail.com") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
